$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- Title ---
Replace-Text "Unraveling the Enigma of the Universe: A Journey into Cosmic Mysteries" `
             "The Art of Governance: Navigating the Labyrinth of Public Administration"

# --- Author name ---
Replace-Text "Elizabeth Anderson" "Ms. Jane Carter"

# --- Author email (kept the two literal "." runs and the "org" run untouched) ---
Replace-Text "elizabeth" "janecarter@educationalhub"
Replace-Text "anderson@astronomicalsociety" "ac"

# --- Body paragraph, first line-group ---
Replace-Text "The boundless expanse of the universe, with its celestial wonders and enigmatic phenomena, has long captivated the minds of humankind" `
             "Embarking on a journey through the intricate world of governance is akin to traversing a labyrinth, where countless paths diverge and converge, each leading to a unique tapestry of societal outcomes"

Replace-Text " From ancient civilizations marveling at the night sky to modern-day astrophysicists probing the depths of space, the quest to unravel the mysteries of the cosmos has fueled scientific inquiry throughout the ages" `
             " Governance, in essence, is the art of steering the ship of state, ensuring its smooth and efficient operation while safeguarding the interests of its citizens"

Replace-Text " This intellectual odyssey has painted the canvas of history with tales of curious observations, profound theories, groundbreaking discoveries, and unsolved riddles. As we continue to venture into the unknown, the exploration of the universe remains a testament to humanity's unwavering determination to comprehend the vastness and complexities of existence" `
             " In this discourse, we shall delve into the complexities of governance, unveiling the profound impact it has on our collective existence"

# --- Body paragraph, second line-group (after first <w:br/>) ---
Replace-Text "The tapestry of the universe, woven with celestial bodies, cosmic phenomena, and fundamental forces, presents an awe-inspiring spectacle for contemplation" `
             "Like a conductor orchestrating a symphony, governance harmonizes the diverse elements of society, ensuring their seamless collaboration towards shared objectives"

Replace-Text " Stars, the luminous beacons of energy, have birthed and nurtured life across the eons, while planets orbit in delicate harmony" `
             " Whether it be the allocation of resources, the enforcement of laws, or the provision of essential services, governance acts as the invisible hand that shapes our communities"

Replace-Text " Galaxies, vast cosmic cities shimmering with countless stars, span light-years, their gravitational embrace shaping the structures of the universe. Yet, amidst this grand spectacle, mysteries persist. Dark matter and dark energy, enigmatic entities that permeate the universe, remain largely incomprehensible, challenging our understanding of gravity and the very nature of space-time" `
             " From the bustling streets of metropolises to the serene landscapes of rural villages, the effects of governance are ubiquitous, affecting every aspect of our lives"

# --- Body paragraph, third line-group (after second <w:br/>) ---
Replace-Text "The ceaseless symphony of the cosmos resounds with cosmic enigmas" `
             "Governance, however, is not a monolithic entity; it manifests itself in myriad forms across different societies"

Replace-Text " Black holes, celestial maelstroms of immense gravitational pull, defy our understanding of time and space" `
             " From the classical Athenian democracy, where citizens directly participated in decision-making, to the modern-day representative democracies, where elected officials serve as the voice of the people, governance has evolved alongside human civilization"

Replace-Text " Supernovas, cataclysmic explosions of dying stars, unveil the raw power of the universe while releasing elements essential for life. Gravitational waves, ripples in the fabric of space-time, provide a glimpse into the universe's most violent events. These celestial phenomena, observed through telescopes and studied with advanced theories, present intriguing puzzles that beckon us to expand the boundaries of our knowledge and push the frontiers of scientific understanding" `
             " Each model bears its own strengths and challenges, reflecting the unique cultural, historical, and socio-economic contexts in which it operates"

# --- Summary paragraph ---
Replace-Text "Unraveling the cosmic mysteries requires the harmonization of scientific inquiry, philosophical contemplation, and artistic expression" `
             "In conclusion, governance stands as the cornerstone of human civilization, providing the framework for peaceful coexistence and collective progress"

Replace-Text " By weaving together observations, theories, and human imagination, we embark on a profound journey to comprehend the universe's grand design" `
             " Its intricate mechanisms, like the gears of a finely tuned machine, orchestrate the complex interactions within society, ensuring order, stability, and the pursuit of common goals"

Replace-Text " As we delve deeper into the enigmatic tapestry of the cosmos, we not only expand our knowledge but also foster a sense of wonder and humility, reminding us of our place amidst the vastness of existence" `
             " Through the examination of governance, we gain a deeper appreciation for the challenges and opportunities inherent in the art of public administration. As citizens, it is our responsibility to actively participate in shaping our governance systems, holding our leaders accountable and working together to build a just, equitable, and sustainable society for generations to come"

# --- Add a new empty paragraph at the very end of the document body ---
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.Text = "`r"
